# Auto-generated edit script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1443227.8
$ws.Range("I9").Value = 1855480.1
$ws.Range("J9").Value = 344.5
$ws.Range("K9").Value = 1855480.1
$ws.Range("L9").Value = 344.5
$ws.Range("M9").Value = -1855311.1
$ws.Range("N9").Value = -682.5
$ws.Range("H11").Value = 691.6667
$ws.Range("I11").Value = 691.6667
$ws.Range("K11").Value = 691.6667
$ws.Range("M11").Value = -551.6667
$ws.Range("H52").Value = 2703
$ws.Range("I52").Value = 4266.5
$ws.Range("J52").Value = 2077.6
$ws.Range("K52").Value = 12799.5
$ws.Range("L52").Value = 6232.799999999999
$ws.Range("M52").Value = -12639.5
$ws.Range("N52").Value = -6552.799999999999
$ws.Range("H70").Value = 102388.7
$ws.Range("I70").Value = 2574.6667
$ws.Range("K70").Value = 7724.000100000001
$ws.Range("M70").Value = -7454.000100000001
$ws.Range("H73").Value = 102388.7
$ws.Range("I73").Value = 2574.6667
$ws.Range("K73").Value = 7724.000100000001
$ws.Range("M73").Value = -6788.000100000001
$ws.Range("H94").Value = 538.6667
$ws.Range("I94").Value = 556.4
$ws.Range("K94").Value = 556.4
$ws.Range("M94").Value = -105.4
$ws.Range("H124").Value = 89999
$ws.Range("J124").Value = 89999
$ws.Range("L124").Value = 89999
$ws.Range("N124").Value = -99819
$ws.Range("H138").Value = 3325.4043
$ws.Range("J138").Value = 4618.1113
$ws.Range("L138").Value = 13854.3339
$ws.Range("N138").Value = -24134.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1663.8276
$ws.Range("I97").Value = 577.1
$ws.Range("K97").Value = 577.1
$ws.Range("M97").Value = -81.10000000000002
$ws.Range("H122").Value = 1633.6666
$ws.Range("I122").Value = 1515.125
$ws.Range("K122").Value = 4545.375
$ws.Range("M122").Value = -2095.375
$ws.Range("H132").Value = 4269.4644
$ws.Range("I132").Value = 4349.5835
$ws.Range("K132").Value = 13048.7505
$ws.Range("M132").Value = -10518.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1265.3914
$ws.Range("I80").Value = 1457.5714
$ws.Range("K80").Value = 1457.5714
$ws.Range("M80").Value = -459.5714
$ws.Range("H83").Value = 1265.3914
$ws.Range("I83").Value = 1457.5714
$ws.Range("K83").Value = 7287.857
$ws.Range("M83").Value = -2295.857
$ws.Range("H86").Value = 5229.5
$ws.Range("I86").Value = 4860.2
$ws.Range("J86").Value = 5493.2856
$ws.Range("K86").Value = 4860.2
$ws.Range("L86").Value = 5493.2856
$ws.Range("M86").Value = -3737.2
$ws.Range("N86").Value = -7739.2856
$ws.Range("H89").Value = 5229.5
$ws.Range("I89").Value = 4860.2
$ws.Range("J89").Value = 5493.2856
$ws.Range("K89").Value = 24301
$ws.Range("L89").Value = 27466.428
$ws.Range("M89").Value = -18685
$ws.Range("N89").Value = -38698.428
$ws.Range("H94").Value = 1075.0834
$ws.Range("I94").Value = 1075.0834
$ws.Range("K94").Value = 1075.0834
$ws.Range("M94").Value = -624.0834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 238
$ws.Range("I7").Value = 200.5
$ws.Range("J7").Value = 263
$ws.Range("K7").Value = 200.5
$ws.Range("L7").Value = 263
$ws.Range("M7").Value = -87.5
$ws.Range("N7").Value = -489
$ws.Range("H16").Value = 24092.166
$ws.Range("I16").Value = 910.8
$ws.Range("K16").Value = 910.8
$ws.Range("M16").Value = -623.8
$ws.Range("H58").Value = 11303.667
$ws.Range("I58").Value = 15955.5
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 15955.5
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -15752.5
$ws.Range("N58").Value = -2406
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 24092.166
$ws.Range("I113").Value = 910.8
$ws.Range("K113").Value = 910.8
$ws.Range("M113").Value = 1259.2
$ws.Range("H122").Value = 431398.16
$ws.Range("I122").Value = 1002319.7
$ws.Range("J122").Value = 3207
$ws.Range("K122").Value = 3006959.1
$ws.Range("L122").Value = 9621
$ws.Range("M122").Value = -3004509.1
$ws.Range("N122").Value = -14521
$ws.Range("H136").Value = 11303.667
$ws.Range("I136").Value = 15955.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 47866.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -45316.5
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1800.3334
$ws.Range("J5").Value = 887.5
$ws.Range("L5").Value = 2662.5
$ws.Range("N5").Value = -2886.5
$ws.Range("H12").Value = 256.13043
$ws.Range("J12").Value = 272.4375
$ws.Range("L12").Value = 817.3125
$ws.Range("N12").Value = -1163.3125
$ws.Range("H92").Value = 449
$ws.Range("I92").Value = 399
$ws.Range("K92").Value = 1197
$ws.Range("M92").Value = 51
$ws.Range("H132").Value = 1314
$ws.Range("I132").Value = 459.5
$ws.Range("J132").Value = 2168.5
$ws.Range("K132").Value = 4135.5
$ws.Range("L132").Value = 19516.5
$ws.Range("M132").Value = -1605.5
$ws.Range("N132").Value = -24576.5
$ws.Range("H135").Value = 1800.3334
$ws.Range("J135").Value = 887.5
$ws.Range("L135").Value = 7987.5
$ws.Range("N135").Value = -13057.5
$ws.Range("H141").Value = 5345.2
$ws.Range("I141").Value = 4181.75
$ws.Range("K141").Value = 12545.25
$ws.Range("M141").Value = -7365.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1477
$ws.Range("I16").Value = 1939.8
$ws.Range("K16").Value = 1939.8
$ws.Range("M16").Value = -1769.8
$ws.Range("H40").Value = 12423.192
$ws.Range("I40").Value = 10252.111
$ws.Range("K40").Value = 10252.111
$ws.Range("M40").Value = -10116.111
$ws.Range("H46").Value = 5238
$ws.Range("J46").Value = 3894
$ws.Range("L46").Value = 3894
$ws.Range("N46").Value = -4270
$ws.Range("H55").Value = 420.36365
$ws.Range("I55").Value = 468.42856
$ws.Range("J55").Value = 336.25
$ws.Range("K55").Value = 468.42856
$ws.Range("L55").Value = 336.25
$ws.Range("M55").Value = -295.42856
$ws.Range("N55").Value = -682.25
$ws.Range("H100").Value = 4718.091
$ws.Range("I100").Value = 2474.75
$ws.Range("K100").Value = 2474.75
$ws.Range("M100").Value = -1933.75
$ws.Range("H132").Value = 3400.5757
$ws.Range("I132").Value = 3316.7856
$ws.Range("K132").Value = 9950.356800000001
$ws.Range("M132").Value = -7420.356800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5823.5713
$ws.Range("I62").Value = 5823.5713
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5823.5713
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5199.5713
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 5823.5713
$ws.Range("I65").Value = 5823.5713
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 29117.8565
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -25997.8565
$ws.Range("N65").ClearContents()
$ws.Range("H96").Value = 166667820
$ws.Range("I96").Value = 200001070
$ws.Range("K96").Value = 200001070
$ws.Range("M96").Value = -199999697
$ws.Range("H98").Value = 15000
$ws.Range("J98").Value = 15000
$ws.Range("L98").Value = 15000
$ws.Range("N98").Value = -20990
$ws.Range("H107").Value = 1791.6666
$ws.Range("I107").Value = 1536.6154
$ws.Range("J107").Value = 2093.0908
$ws.Range("K107").Value = 4609.8462
$ws.Range("L107").Value = 6279.2724
$ws.Range("M107").Value = -2689.8462
$ws.Range("N107").Value = -10119.2724
$ws.Range("H122").Value = 4935.65
$ws.Range("I122").Value = 4372.9443
$ws.Range("K122").Value = 13118.8329
$ws.Range("M122").Value = -10668.8329
